# Update the title / body / source-link text on slides 2-11 of the
# Python presentation (slide indices are 1-based; slide 1 is the title
# slide and is left untouched).

$p = $ppt.ActivePresentation

function Set-SlideText($Slide, $Title, $Body, $Link) {
    $Slide.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = $Title
    $Slide.Shapes.Item(2).TextFrame.TextRange.Runs(1,1).Text = $Body

    $linkShape = $Slide.Shapes.Item(3)
    $originalHeight = $linkShape.Height
    $linkShape.TextFrame.TextRange.Paragraphs(2,1).Runs(1,1).Text = $Link
    # Editing the text re-triggers the textbox's auto-fit-to-text sizing;
    # restore the original (unchanged, per the source diff) box height.
    $linkShape.Height = $originalHeight
}

# Slide 2: Introduction to Python
Set-SlideText $p.Slides.Item(2) `
    "Introduction to Python" `
    "Python is a high-level programming language known for its simple syntax and readability. It is widely used in various applications including web development, data analysis, artificial intelligence, and more." `
    "- Python official website: www.python.org"

# Slide 3: Python Features
Set-SlideText $p.Slides.Item(3) `
    "Python Features" `
    "Python has numerous features such as dynamic typing, automatic memory management, extensive standard library, and more. These features make Python a versatile and powerful language." `
    "- Python documentation: docs.python.org"

# Slide 4: Python Data Types
Set-SlideText $p.Slides.Item(4) `
    "Python Data Types" `
    "Python supports various data types including integers, floats, strings, lists, tuples, dictionaries, and more. Understanding data types is crucial for writing efficient and effective Python code." `
    "- Python data types documentation: docs.python.org/library/stdtypes.html"

# Slide 5: Control Structures in Python
Set-SlideText $p.Slides.Item(5) `
    "Control Structures in Python" `
    "Python provides control structures such as loops (for, while) and conditional statements (if, else, elif) for controlling the flow of a program. These structures help in making decisions and repeating tasks." `
    "- Python control structures documentation: docs.python.org/tutorial/controlflow.html"

# Slide 6: Functions in Python
Set-SlideText $p.Slides.Item(6) `
    "Functions in Python" `
    "Functions in Python allow us to encapsulate code for reusability and modularity. They help in organizing code and making it easier to manage and maintain." `
    "- Python functions documentation: docs.python.org/tutorial/controlflow.html#defining-functions"

# Slide 7: Python Libraries
Set-SlideText $p.Slides.Item(7) `
    "Python Libraries" `
    "Python has a vast collection of libraries such as NumPy, Pandas, Matplotlib, and TensorFlow that extend its functionality for specific tasks. These libraries make Python a popular choice for data analysis, machine learning, and more." `
    "- Python libraries documentation: numpy.org, pandas.pydata.org, matplotlib.org, tensorflow.org"

# Slide 8: Object-Oriented Programming in Python
Set-SlideText $p.Slides.Item(8) `
    "Object-Oriented Programming in Python" `
    "Python supports object-oriented programming principles such as encapsulation, inheritance, and polymorphism. Classes and objects are fundamental concepts in Python for building reusable and modular code." `
    "- Python classes documentation: docs.python.org/tutorial/classes.html"

# Slide 9: Python Development Environments
Set-SlideText $p.Slides.Item(9) `
    "Python Development Environments" `
    "There are various IDEs (Integrated Development Environments) and text editors available for Python development such as PyCharm, Visual Studio Code, and Jupyter Notebook. Choosing the right environment can enhance productivity and efficiency in coding." `
    "- PyCharm: jetbrains.com/pycharm, Visual Studio Code: code.visualstudio.com, Jupyter Notebook: jupyter.org"

# Slide 10: Python Community and Resources
Set-SlideText $p.Slides.Item(10) `
    "Python Community and Resources" `
    "Python has a vibrant and supportive community with numerous forums, online resources, and tutorials available for help and learning. Engaging with the Python community can help in networking, learning new concepts, and solving coding challenges." `
    "- Python community website: python.org/community, Python subreddit: reddit.com/r/python"

# Slide 11: Conclusion
Set-SlideText $p.Slides.Item(11) `
    "Conclusion" `
    "Python is a versatile and powerful programming language with a wide range of applications. Its simplicity, readability, and extensive library support make it a popular choice for developers worldwide." `
    "- Python summary: www.python.org/about"
